$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: input_variables changed from "PAL" to "pal"
$ws.Range("F8").Value = "pal"

# Rows 97-108: fill in missing input_variables (column F) values to mirror
# the dataschema_variable (column B) for each row.
$ws.Range("F97").Value = "SUGAR_CONFECT_11"
$ws.Range("F98").Value = "CAKES_12"
$ws.Range("F99").Value = "FRUITVEG_JUICE_1301"
$ws.Range("F100").Value = "SOFTDRINKS_1302"
$ws.Range("F101").Value = "ART_SWEETENER_170201"
$ws.Range("F102").Value = "VEGETABLES_02"
$ws.Range("F103").Value = "LEGUMES_TOT_03"
$ws.Range("F104").Value = "FRUITS_TOT_04"
$ws.Range("F105").Value = "RED_MEAT_0701"
$ws.Range("F106").Value = "PROCMEAT_0704"
$ws.Range("F107").Value = "COFFEE_130301"
$ws.Range("F108").Value = "TEA_130302"
